# Reorder worksheet tabs: move "总计" (the summary sheet, currently 2nd)
# in front of "2021-Q3" (currently 1st), so "总计" becomes the first /
# active tab and "2021-Q3" becomes the second tab. No cell data changes.

$wb = $excel.ActiveWorkbook

$summarySheet = $wb.Worksheets.Item("总计")
$firstSheet   = $wb.Worksheets.Item(1)

$summarySheet.Move($firstSheet)
